$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "description"

# --- id column (A2:A6) ---
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5

# --- description column (B2:B6), English values in new order ---
$ws.Range("B2").Value = "studying"
$ws.Range("B3").Value = "reserved"
$ws.Range("B4").Value = "suspended"
$ws.Range("B5").Value = "dropout"
$ws.Range("B6").Value = "graduated"

# --- SQL insert-statement helper formulas (column C) ---
$ws.Range("C2").Formula = "=CONCAT(""INSERT INTO status (description) VALUE ('"",B2,""');"")"
$ws.Range("C3:C6").Formula = "=CONCAT(""INSERT INTO status (description) VALUE ('"",B3,""');"")"

# --- Formatting ---
# Header row + id column: left/center alignment
$ws.Range("A1:B1").HorizontalAlignment = -4131
$ws.Range("A1:B1").VerticalAlignment = -4108
$ws.Range("A2:A6").HorizontalAlignment = -4131
$ws.Range("A2:A6").VerticalAlignment = -4108

# description + helper column: Arial 10, left/center alignment
$ws.Range("B2:B6").Font.Name = "Arial"
$ws.Range("B2:B6").Font.Size = 10
$ws.Range("B2:B6").Font.Color = 0
$ws.Range("B2:B6").HorizontalAlignment = -4131
$ws.Range("B2:B6").VerticalAlignment = -4108

# --- Empty placeholder column D, styled like column B ---
$ws.Range("D1").HorizontalAlignment = -4131
$ws.Range("D1").VerticalAlignment = -4108

$ws.Range("D2:D6").Font.Name = "Arial"
$ws.Range("D2:D6").Font.Size = 10
$ws.Range("D2:D6").Font.Color = 0
$ws.Range("D2:D6").HorizontalAlignment = -4131
$ws.Range("D2:D6").VerticalAlignment = -4108

# --- Column widths ---
$ws.Columns(2).ColumnWidth = 14.0
$ws.Columns(3).ColumnWidth = 49.15
$ws.Columns(4).ColumnWidth = 14.0

# --- Page setup ---
$ws.PageSetup.Orientation = 1

# --- Selection ---
$ws.Range("F8").Select()

Write-Host "Edit applied"
